$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(52)

# Rename UndoRedoStack -> UndoRedoCareTaker
$sh.TextFrame.TextRange.Text = "UndoRedoCareTaker"

# Widen/reposition the shape to fit the new label (left shifts left, width grows; top/height unchanged)
$sh.Left = 96.0
$sh.Width = 108.0001
